$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string for the navigation menu task description
$newText = "A képek és ikonok letöltése és rendszerezése."

# Row 7 previously had A7 numbered "5" but empty B7; the new content
# pushes the numbering down: A7 becomes 5 (style like A5/A6 group: style index 5)
$ws.Range("A7").Style = $ws.Range("A6").Style
$ws.Range("B7").Style = $ws.Range("B6").Style
$ws.Range("B7").Value = $newText

# B5 should use the bordered style (same cellXf family as other data rows, index 3)
$ws.Range("B5").Style = $ws.Range("B4").Style

# Update the active selection to B8 (one row further down) to mirror the diff
$ws.Range("B8").Select()
